$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 1917
$ws.Range("K3").Value = 1824
$ws.Range("B4").Value = 1697
$ws.Range("K4").Value = 394
$ws.Range("K5").Value = 120
$ws.Range("K6").Value = 2351
$ws.Range("B7").Value = 23330
$ws.Range("K7").Value = 6606

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 127
$ws.Range("K3").Value = 123
$ws.Range("K4").Value = 23
$ws.Range("K6").Value = 151
$ws.Range("K7").Value = 434

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 58
$ws.Range("K3").Value = 44
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 143

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 102
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 268

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 55
$ws.Range("K3").Value = 72
$ws.Range("K6").Value = 74
$ws.Range("K7").Value = 220

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 45
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 160

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 52
$ws.Range("K6").Value = 52
$ws.Range("K7").Value = 183
$ws.Range("K8").Value = 434
$ws.Range("K11").Value = 138
$ws.Range("K15").Value = 61
$ws.Range("K19").Value = 184
$ws.Range("K20").Value = 144
$ws.Range("K22").Value = 18
$ws.Range("K23").Value = 60
$ws.Range("K27").Value = 73
$ws.Range("K29").Value = 325
$ws.Range("K31").Value = 75
$ws.Range("K33").Value = 268
$ws.Range("K37").Value = 220
$ws.Range("K42").Value = 230
$ws.Range("K43").Value = 64
$ws.Range("K44").Value = 65
$ws.Range("K48").Value = 78
$ws.Range("K52").Value = 179
$ws.Range("K54").Value = 111
$ws.Range("K60").Value = 45
$ws.Range("B63").Value = 402
$ws.Range("K63").Value = 23
$ws.Range("K64").Value = 44
$ws.Range("K65").Value = 160
$ws.Range("K67").Value = 256
$ws.Range("K72").Value = 27
$ws.Range("K75").Value = 26
$ws.Range("K76").Value = 96
$ws.Range("K77").Value = 44
$ws.Range("K78").Value = 85
$ws.Range("K79").Value = 176
$ws.Range("K83").Value = 143
$ws.Range("K85").Value = 331
$ws.Range("K88").Value = 87
$ws.Range("K90").Value = 57
$ws.Range("K93").Value = 29
$ws.Range("K96").Value = 93
$ws.Range("K97").Value = 60
$ws.Range("K98").Value = 45
$ws.Range("B101").Value = 23330
$ws.Range("K101").Value = 6606

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 27
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 73
$ws.Range("K3").Value = 80
$ws.Range("K6").Value = 84
$ws.Range("K7").Value = 256

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 38
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 85
$ws.Range("K3").Value = 109
$ws.Range("K6").Value = 108
$ws.Range("K7").Value = 325

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K3").Value = 13
$ws.Range("K4").Value = 13
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 54
$ws.Range("K6").Value = 61
$ws.Range("K7").Value = 184

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K2").Value = 11
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 65

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K2").Value = 18
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 96

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K2").Value = 19
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K3").Value = 62
$ws.Range("K6").Value = 101
$ws.Range("K7").Value = 230

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 85

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K3").Value = 61
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 43
$ws.Range("K7").Value = 144

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 29

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 56
$ws.Range("K7").Value = 183

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 19
$ws.Range("K7").Value = 61

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 44
$ws.Range("K7").Value = 138

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 17
$ws.Range("K6").Value = 49
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 26

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 26
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 120
$ws.Range("K7").Value = 331

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 18

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 45
$ws.Range("K3").Value = 39
$ws.Range("K7").Value = 179

Write-Host "Applied 152 cell updates across 41 sheets"
